# Updated output and input excel
#
# Applies the wishing-list workbook update:
#   1. Appends "Alina Tamminen" to Julius Kuusisto's wish list (row 134, col C).
#   2. Appends "Julius Kuusisto" to Alina Tamminen's wish list (row 137, col C).
#   3. Adds four new sample rows (a/b/c/d wishing-list chain) after the
#      existing data (rows 143-146).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend Julius Kuusisto's wish-list entry ---
$ws.Range("C134").Value = "Kim Kuusisto, Kati Kuusisto, Margareta Kuusisto, Topi Kuusisto, Tanja Laurila, Alina Tamminen"

# --- 2. Extend Alina Tamminen's wish-list entry ---
$ws.Range("C137").Value = "Johan Tamminen, Elli Tamminen, Konsta Tamminen, Marja-Liisa Tamminen, Topi Kuusisto, Julius Kuusisto"

# --- 3. Append the new sample participant/wishing-list rows ---
$ws.Range("B143").Value = "a"
$ws.Range("C143").Value = "b"

$ws.Range("B144").Value = "b"
$ws.Range("C144").Value = "c"

$ws.Range("B145").Value = "c"
$ws.Range("C145").Value = "b, d"

$ws.Range("B146").Value = "d"
$ws.Range("C146").Value = "a, c"

# --- Restore the selected cell to match the author's saved cursor position ---
$ws.Range("C132").Select() | Out-Null
